$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp on the existing last row (row 49, column A) with
# its re-retrieved (slightly different) fractional-day value.
$ws.Cells.Item(49, 1).Value = 44362.76909272106

# Append the newly retrieved data row as row 50.
$ws.Cells.Item(50, 1).Value = 44363.77004171583
$ws.Cells.Item(50, 2).Value = 78178
$ws.Cells.Item(50, 3).Value = 65747
$ws.Cells.Item(50, 4).Value = 3510
$ws.Cells.Item(50, 5).Value = 2109
$ws.Cells.Item(50, 6).Value = 1490
$ws.Cells.Item(50, 7).Value = 20631
$ws.Cells.Item(50, 8).Value = 1498
$ws.Cells.Item(50, 9).Value = 897
$ws.Cells.Item(50, 10).Value = 191
